# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same shape as the
# existing quarterly sheets) right before the "总计" (grand-total) sheet, and
# refreshes "总计" with a new leading row summarizing 2022-Q1.

function Set-HeaderCell($cell, $value) {
    # Bold, centered, thin-bordered style used for header row / index column
    # on the quarterly detail + total sheets.
    $cell.Value = $value
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

function Set-TextCell($cell, $value) {
    # Force text storage so numeric-looking strings (fund codes, percentages)
    # keep leading zeros / exact formatting instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted right before "总计".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

Set-HeaderCell $q1.Cells.Item(1, 2) "基金代码"
Set-HeaderCell $q1.Cells.Item(1, 3) "基金名称"
Set-HeaderCell $q1.Cells.Item(1, 4) "基金规模"
Set-HeaderCell $q1.Cells.Item(1, 5) "股票总仓位"
Set-HeaderCell $q1.Cells.Item(1, 6) "仓位占比"
Set-HeaderCell $q1.Cells.Item(1, 7) "持有市值(亿元)"
Set-HeaderCell $q1.Cells.Item(1, 8) "仓位排名"

$q1Data = @(
    @(0, "011284", "中信保诚龙腾精选混合",       "1.22", "75.38", "2.64", "0.0322", 8),
    @(1, "006209", "中信保诚新蓝筹灵活配置混合", "1.16", "77.03", "2.67", "0.0310", 8),
    @(2, "011603", "兴业高端制造混合A",           "1.19", "76.40", "2.54", "0.0302", 10),
    @(3, "011604", "兴业高端制造混合C",           "0.54", "76.40", "2.54", "0.0137", 10),
    @(4, "002453", "九泰久稳灵活配置混合A",       "0.09", "94.85", "3.55", "0.0032", 2),
    @(5, "002454", "九泰久稳灵活配置混合C",       "0.04", "94.85", "3.55", "0.0014", 2)
)

$row = 2
foreach ($item in $q1Data) {
    Set-HeaderCell $q1.Cells.Item($row, 1) $item[0]
    Set-TextCell $q1.Cells.Item($row, 2) $item[1]
    Set-TextCell $q1.Cells.Item($row, 3) $item[2]
    Set-TextCell $q1.Cells.Item($row, 4) $item[3]
    Set-TextCell $q1.Cells.Item($row, 5) $item[4]
    Set-TextCell $q1.Cells.Item($row, 6) $item[5]
    Set-TextCell $q1.Cells.Item($row, 7) $item[6]
    $q1.Cells.Item($row, 8).Value = $item[7]
    $row++
}

# ---------------------------------------------------------------------
# 2. Refresh "总计" with a new leading 2022-Q1 row (existing rows shift
#    down, index column renumbers 0..3).
# ---------------------------------------------------------------------
# NOTE: re-fetch "总计" fresh here (rather than reusing a reference
# captured before the Worksheets.Add above) -- the COM shim repoints a
# captured worksheet handle to track the newly-inserted sheet once
# Add()/rename happens, so an earlier-bound variable would silently
# write into "2022-Q1" instead.
$totalSheet = $wb.Worksheets.Item("总计")

$totalData = @(
    @(0, "2022-Q1", 6, 0.11),
    @(1, "2021-Q4", 14, 2.76),
    @(2, "2021-Q3", 12, 2.05),
    @(3, "2021-Q2", 2, 0.05)
)

$row = 2
foreach ($item in $totalData) {
    $totalSheet.Cells.Item($row, 1).Value = $item[0]
    $totalSheet.Cells.Item($row, 2).Value = $item[1]
    $totalSheet.Cells.Item($row, 3).Value = $item[2]
    $totalSheet.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# Row 5 (A5) is brand new territory for this sheet -- give it the same
# bold/border/center-top styling as the pre-existing index cells A2:A4.
$a5 = $totalSheet.Cells.Item(5, 1)
$a5.Font.Bold = $true
$a5.HorizontalAlignment = -4108
$a5.VerticalAlignment = -4160
$a5.Borders.Item(7).LineStyle = 1
$a5.Borders.Item(8).LineStyle = 1
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(10).LineStyle = 1

# Restore the originally-active sheet/tab (inserting a sheet shifts focus
# onto it by default).
$wb.Worksheets.Item("2021-Q2").Activate()
